# Actualización de datos obtenidos el 6 de abril de 2016
# "Grado" and "Sexo" were previously classified as measures; they are now
# reclassified as dimensions, with accompanying code-list mapping files.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D = Grado
$ws.Range("D3").Value = "iaest-dimension:grado"   # was iaest-measure:grado
$ws.Range("D4").Value = "dim"                      # was medida
$ws.Range("D5").Value = "skos:Concept"              # was xsd:string
$ws.Range("D6").Value = "mapping-grado.xlsx"        # new row

# Column F = Sexo
$ws.Range("F3").Value = "iaest-dimension:sexo"     # was iaest-measure:sexo
$ws.Range("F4").Value = "dim"                       # was medida
$ws.Range("F5").Value = "skos:Concept"               # was xsd:string
$ws.Range("F6").Value = "mapping-sexo.xlsx"          # new row
